$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Source")

# Add a new row with a single-space value in column B (check for empty import)
$ws.Range("B4").Value = " "

# Update the selection to match the post-edit state
$ws.Range("E12").Select()
